# Add a new "ATS History" worksheet (Grok-verified) as the last sheet in the
# workbook, populated with the ATS adoption history table + key insight.

$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore the selection
# after appending the new tab at the end of the workbook.
$firstSheet = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ATS History"

# Match the page margins used throughout the rest of the workbook
# (0.75in left/right, 1in top/bottom, 0.5in header/footer).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Title
$ws.Range("A1").Value = "ATS ADOPTION HISTORY (Grok-verified)"

# Table header (row 3 - row 2 left blank as a spacer like the other sheets)
$ws.Range("A3").Value = "Era"
$ws.Range("B3").Value = "Milestone"
$ws.Range("C3").Value = "Source"

# Table rows
$ws.Range("A4").Value = "Late 1970s-80s"
$ws.Range("B4").Value = "Early HR databases on mainframes (IBM)"
$ws.Range("C4").Value = "SHRM"

# Era values that look numeric ("1990", "1999", "2002") must stay plain
# text, matching the source file, rather than being auto-coerced into
# numbers. Prefix with an apostrophe to force text entry, then reset the
# cell style back to Normal so no stray "quote prefix" formatting sticks.
$ws.Range("A5").Value = "'1990"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "First commercial ATS - Resumix founded"
$ws.Range("C5").Value = "Harvard Business Review"

$ws.Range("A6").Value = "'1999"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "Taleo founded (web-based ATS)"
$ws.Range("C6").Value = "Oracle"

$ws.Range("A7").Value = "'2002"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "50%+ of large US companies using ATS"
$ws.Range("C7").Value = "SHRM"

$ws.Range("A8").Value = "2010s"
$ws.Range("B8").Value = "AI-driven screening, mobile integration"
$ws.Range("C8").Value = "Industry reports"

$ws.Range("A9").Value = "Today"
$ws.Range("B9").Value = "99% of Fortune 500 use ATS"
$ws.Range("C9").Value = "Indeed Career Guide"

# Key insight (row 10 left blank as a spacer)
$ws.Range("A11").Value = "KEY INSIGHT"
$ws.Range("A12").Value = "ATS has been filtering resumes for 35+ years. MatchForge checks against 10 systems (~56% market)."

# Restore the original active sheet/selection.
$firstSheet.Activate()
